# Edit script: Add/update resource data for UWTranslationQuestions
$d = $word.ActiveDocument

# 1. Delete the "License Information" Heading2 paragraph entirely.
$licensePara = $d.Paragraphs.Item(4)
if ($licensePara.Range.Text -like "License Information*") {
    $licensePara.Range.Delete()
}

# 2. Replace the body of the "is based on..." paragraph (now paragraph 4)
#    with the new, simplified license/adaptation text, preserving the
#    paragraph's own formatting (w:bidi) by including it in the replacement XML.
$target = $d.Paragraphs.Item(4)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:bidi/></w:pPr><w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr></w:r><w:r><w:rPr><w:b/><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>unfoldingWord® Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t xml:space="preserve"> © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. </w:t></w:r><w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>unfoldingWord® Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t xml:space="preserve"> has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from </w:t></w:r><w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>unfoldingWord® Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t xml:space="preserve"> © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual</w:t></w:r><w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($newParaXml)

# 3. Delete the "This PDF version is provided under the same license." paragraph.
$pdfPara = $d.Paragraphs.Item(5)
if ($pdfPara.Range.Text -like "This PDF version*") {
    $pdfPara.Range.Delete()
}
